# Generate Report for Handback
#
# This script updates the localization-status workbook so that the
# "fa1d7f6e-7967-403a-896f-c7c2b2290ebc" row reflects a failed handback
# transform instead of "Ready for handoff", and records the handback
# error detail for both the zh-cn and de-de target languages.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: the status shown for zh-cn (E3) and de-de (F3) for the
# fa1d7f6e-... file moves from "Ready for handoff" to "Handback transform failed"
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# Per-language sheets: the "Status" column (C) for the same file row (row 3)
$ws2.Range("C3").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# Error Detail column (P) gets widened to fit the new error messages
$colWidth = 40 - 5/6
$ws2.Range("P1").EntireColumn.ColumnWidth = $colWidth
$ws3.Range("P1").EntireColumn.ColumnWidth = $colWidth

# Error Detail column (P), row 3, now records why the handback failed
$ws2.Range("P3").Value = "Handback file name: hqhpqj15.c5g is different with handoff file name: fa1d7f6e-7967-403a-896f-c7c2b2290ebc.ab6741b01f672bc0c5f65d5c114e35e18de6bc6f.zh-cn."
$ws3.Range("P3").Value = "Handback file name: hqhpqj15.c5g is different with handoff file name: fa1d7f6e-7967-403a-896f-c7c2b2290ebc.ab6741b01f672bc0c5f65d5c114e35e18de6bc6f.de-de."
